$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 34486036
$ws.Range("I64").Value = 62502564
$ws.Range("J64").Value = 4158.4614
$ws.Range("K64").Value = 62502564
$ws.Range("L64").Value = 4158.4614
$ws.Range("M64").Value = -62502316
$ws.Range("N64").Value = -4654.4614

$ws.Range("H67").Value = 34486036
$ws.Range("I67").Value = 62502564
$ws.Range("J67").Value = 4158.4614
$ws.Range("K67").Value = 62502564
$ws.Range("L67").Value = 4158.4614
$ws.Range("M67").Value = -62501706
$ws.Range("N67").Value = -5874.4614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 295416.2
$ws.Range("I32").Value = 323764.97
$ws.Range("J32").Value = 11928.571
$ws.Range("K32").Value = 323764.97
$ws.Range("L32").Value = 11928.571
$ws.Range("M32").Value = -323477.97
$ws.Range("N32").Value = -12502.571

$ws.Range("H45").Value = 2224.4614
$ws.Range("I45").Value = 1646.2858
$ws.Range("J45").Value = 2899
$ws.Range("K45").Value = 1646.2858
$ws.Range("L45").Value = 2899
$ws.Range("M45").Value = -1269.2858
$ws.Range("N45").Value = -3653

$ws.Range("H139").Value = 47266.668
$ws.Range("J139").Value = 47266.668
$ws.Range("L139").Value = 47266.668
$ws.Range("N139").Value = -57546.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H115").Value = 48000
$ws.Range("J115").Value = 48000
$ws.Range("L115").Value = 48000
$ws.Range("N115").Value = -51134

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 566.6667
$ws.Range("I14").Value = 350
$ws.Range("J14").Value = 1000
$ws.Range("K14").Value = 350
$ws.Range("L14").Value = 1000
$ws.Range("M14").Value = -180
$ws.Range("N14").Value = -1340

$ws.Range("H31").Value = 3436.7144
$ws.Range("I31").Value = 1233.0286
$ws.Range("J31").Value = 4660.984
$ws.Range("K31").Value = 1233.0286
$ws.Range("L31").Value = 4660.984
$ws.Range("M31").Value = -938.0286000000001
$ws.Range("N31").Value = -5250.984

$ws.Range("H34").Value = 3436.7144
$ws.Range("I34").Value = 1233.0286
$ws.Range("J34").Value = 4660.984
$ws.Range("K34").Value = 1233.0286
$ws.Range("L34").Value = 4660.984
$ws.Range("M34").Value = -1031.0286
$ws.Range("N34").Value = -5064.984

$ws.Range("H37").Value = 25000
$ws.Range("J37").Value = 25000
$ws.Range("L37").Value = 25000
$ws.Range("N37").Value = -25214

$ws.Range("H62").Value = 4092.8235
$ws.Range("I62").Value = 4558
$ws.Range("J62").Value = 2800.6667
$ws.Range("K62").Value = 4558
$ws.Range("L62").Value = 2800.6667
$ws.Range("M62").Value = -3934
$ws.Range("N62").Value = -4048.6667

$ws.Range("H65").Value = 4092.8235
$ws.Range("I65").Value = 4558
$ws.Range("J65").Value = 2800.6667
$ws.Range("K65").Value = 22790
$ws.Range("L65").Value = 14003.3335
$ws.Range("M65").Value = -19670
$ws.Range("N65").Value = -20243.3335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 769.64703
$ws.Range("I131").Value = 273.8889
$ws.Range("J131").Value = 948.12
$ws.Range("K131").Value = 821.6667
$ws.Range("L131").Value = 2844.36
$ws.Range("M131").Value = 4218.3333
$ws.Range("N131").Value = -12924.36

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 658.2
$ws.Range("I3").Value = 454
$ws.Range("J3").Value = 1134.6666
$ws.Range("K3").Value = 454
$ws.Range("L3").Value = 1134.6666
$ws.Range("M3").Value = -338
$ws.Range("N3").Value = -1366.6666

$ws.Range("H10").Value = 667666.7
$ws.Range("I10").Value = 1000500
$ws.Range("J10").Value = 2000
$ws.Range("K10").Value = 1000500
$ws.Range("L10").Value = 2000
$ws.Range("M10").Value = -1000331
$ws.Range("N10").Value = -2338

$ws.Range("H40").Value = 2500
$ws.Range("I40").Value = 2500
$ws.Range("K40").Value = 2500
$ws.Range("M40").Value = -2349

$ws.Range("H80").Value = 2800.1538
$ws.Range("J80").Value = 3767.6
$ws.Range("L80").Value = 3767.6
$ws.Range("N80").Value = -5763.6

$ws.Range("H83").Value = 2800.1538
$ws.Range("J83").Value = 3767.6
$ws.Range("L83").Value = 18838
$ws.Range("N83").Value = -28822

$ws.Range("H122").Value = 2731.4375
$ws.Range("I122").Value = 1579.4
$ws.Range("J122").Value = 4651.5
$ws.Range("K122").Value = 4738.200000000001
$ws.Range("L122").Value = 13954.5
$ws.Range("M122").Value = -2288.200000000001
$ws.Range("N122").Value = -18854.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 462.9091
$ws.Range("I9").Value = 342.22223
$ws.Range("J9").Value = 1006
$ws.Range("K9").Value = 342.22223
$ws.Range("L9").Value = 1006
$ws.Range("M9").Value = -118.22223
$ws.Range("N9").Value = -1454

$ws.Range("H11").Value = 45003.5
$ws.Range("J11").Value = 45003.5
$ws.Range("L11").Value = 45003.5
$ws.Range("N11").Value = -45283.5

$ws.Range("H22").Value = 20680.4
$ws.Range("I22").Value = 600
$ws.Range("J22").Value = 25700.5
$ws.Range("K22").Value = 600
$ws.Range("L22").Value = 25700.5
$ws.Range("M22").Value = -305
$ws.Range("N22").Value = -26290.5

$ws.Range("H27").Value = 20680.4
$ws.Range("I27").Value = 600
$ws.Range("J27").Value = 25700.5
$ws.Range("K27").Value = 600
$ws.Range("L27").Value = 25700.5
$ws.Range("M27").Value = -493
$ws.Range("N27").Value = -25914.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 3000
$ws.Range("J18").Value = 3000
$ws.Range("L18").Value = 3000
$ws.Range("N18").Value = -3346

$ws.Range("H39").Value = 6445
$ws.Range("J39").Value = 6445
$ws.Range("L39").Value = 6445
$ws.Range("N39").Value = -7271

$ws.Range("H43").Value = 13150
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()

$ws.Range("H122").Value = 13334417
$ws.Range("I122").Value = 20000846
$ws.Range("J122").Value = 1560
$ws.Range("K122").Value = 60002538
$ws.Range("L122").Value = 4680
$ws.Range("M122").Value = -60000088
$ws.Range("N122").Value = -9580
